$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B (session_id) - single UUID shared across all data rows
$ws.Range("B2:B49").Value = "b588e64b-0f4d-4652-b493-1cc58046c7a3"

# Column C (session_entry_id) - two UUID groups
$ws.Range("C2:C22").Value = "e651a4ea-4066-4a6a-af95-bd9fa4572fca"
$ws.Range("C23:C49").Value = "babf7812-8031-4404-bd0c-d5bde182d5a2"

# Column A (ingest_session_issue_id) - unique UUID per row
$ws.Range("A2").Value = "61de82d9-81d2-4e36-ad94-fdba4e09b92c"
$ws.Range("A3").Value = "39fa1eb5-2127-4f2f-ac8a-44d70e3c0643"
$ws.Range("A4").Value = "ced02b99-c1ad-4cc1-8b3d-9b07f30cf603"
$ws.Range("A5").Value = "3bc8b78f-68c6-421f-bd7e-31fddc4f8b25"
$ws.Range("A6").Value = "9571ff4e-1aca-41ff-801d-469a911a081d"
$ws.Range("A7").Value = "31df5ca8-09d9-4144-addb-03ccb645ce53"
$ws.Range("A8").Value = "c7b36a7a-cf9b-4db4-83b1-e7c88648370c"
$ws.Range("A9").Value = "cb98c50f-979b-4a01-b738-2cb41be35bc3"
$ws.Range("A10").Value = "ff433581-43c0-46ca-b9b8-59d11f9b7c55"
$ws.Range("A11").Value = "454ee3e9-bdfb-4a1e-bb75-24c6ff3393ca"
$ws.Range("A12").Value = "bedb2365-801e-4da9-952e-d87801e25365"
$ws.Range("A13").Value = "0dd103ce-8047-4567-b433-3591c9d113f4"
$ws.Range("A14").Value = "72dc07af-0c2e-4c52-9a28-f34cc03b8730"
$ws.Range("A15").Value = "fb0318cb-8bf1-4d08-881f-f2ea63dfee60"
$ws.Range("A16").Value = "ed9c2b9a-de2a-4202-a133-beac12a8d66c"
$ws.Range("A17").Value = "abcb7e43-7726-4000-bd0a-d76330a98b3c"
$ws.Range("A18").Value = "2f945700-6f2d-443c-b9cc-c2ef76d72afd"
$ws.Range("A19").Value = "5f47293f-97d8-4bee-a90e-7eedd5d515bf"
$ws.Range("A20").Value = "d4de0bbf-89cd-480a-b280-be41fe73ca95"
$ws.Range("A21").Value = "9e7173a6-e3b9-49d4-aac0-5dba9d94b68b"
$ws.Range("A22").Value = "efc42a9d-f48a-4c46-9a98-be8df45c43ab"
$ws.Range("A23").Value = "9b4f3391-37bb-48b1-8168-74ae07020e17"
$ws.Range("A24").Value = "62795c5b-957d-4972-a883-7f3bbf0e13cf"
$ws.Range("A25").Value = "27292fa8-ef18-4731-8f42-b07ac50de4da"
$ws.Range("A26").Value = "7ee6ac80-58c9-4227-8a42-7d9e421edfe3"
$ws.Range("A27").Value = "2b04d28c-535c-4ec9-933a-0fa940f38617"
$ws.Range("A28").Value = "73b2fdf4-3a8a-4806-ac5f-748642d3bc42"
$ws.Range("A29").Value = "dc2c90a2-6e33-4684-b18f-befbab41fa32"
$ws.Range("A30").Value = "cd68ccb6-4f04-4215-b304-6ddd5b9afdb4"
$ws.Range("A31").Value = "cc0ed134-9817-4063-8f92-3a18a66ebd4d"
$ws.Range("A32").Value = "88d609c5-18c6-455a-9235-9f7a218f819e"
$ws.Range("A33").Value = "44a91b12-070e-4606-bab8-f8616ebfb984"
$ws.Range("A34").Value = "7c4aea56-6d76-48d3-966a-e4b1f67e50ee"
$ws.Range("A35").Value = "c8a74d18-8263-40d5-b989-cc6fe0e742bf"
$ws.Range("A36").Value = "5b6a6ba1-3358-4621-893e-a05750d9dca1"
$ws.Range("A37").Value = "3a989677-e35b-4af4-99aa-bc0c7d6f7794"
$ws.Range("A38").Value = "683ce513-c08b-4ba4-8e9d-f9b9052efe66"
$ws.Range("A39").Value = "2739a781-d898-4f1c-8c1e-2c178b0baee0"
$ws.Range("A40").Value = "f64ec13d-69be-4f45-8715-c852caa2a558"
$ws.Range("A41").Value = "8448606e-142f-4887-a3cf-0b68c592c029"
$ws.Range("A42").Value = "a3d9d1e4-ece9-43f5-af7e-ef051ee6b808"
$ws.Range("A43").Value = "f8e77cc7-56dd-44de-9fb4-91fa105f213e"
$ws.Range("A44").Value = "eaf463b6-8c8e-4f47-b37e-c6afd3ec1354"
$ws.Range("A45").Value = "1e47eb4f-1e5d-45cf-9eb5-0aef4288e42a"
$ws.Range("A46").Value = "38ba4023-690d-47de-b7b2-1ac31a3ba5a9"
$ws.Range("A47").Value = "4cf31e4f-0233-4459-a7e9-3aad12fda0ca"
$ws.Range("A48").Value = "c32d35c9-2c80-471a-bb2d-e6bfb9ef0e64"
$ws.Range("A49").Value = "e2ec9740-3723-487b-a277-4425efc0aa56"
